$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-11 from 45204 to 45205
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
